$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.071.64"
$ws.Range("E2").Value = "  +3.26%  "

$ws.Range("D3").Value = "3.217.94"
$ws.Range("E3").Value = "  +2.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +6.36%  "

$ws.Range("D9").Value = "3.217.65"
$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("E10").Value = "  +2.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.91%  "

$ws.Range("E12").Value = "  +4.22%  "

$ws.Range("E13").Value = "  +2.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.85%  "

$ws.Range("D15").Value = "3.744.95"
$ws.Range("E15").Value = "  +2.00%  "

$ws.Range("D16").Value = "66.958.17"
$ws.Range("E16").Value = "  +3.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.69%  "

$ws.Range("D18").Value = "3.215.02"
$ws.Range("E18").Value = "  +1.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "528.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.77%  "

$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.748"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.06%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("E28").Value = "  +3.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.62%  "

$ws.Range("E33").Value = "  +3.69%  "

$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.63"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "524.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.47%  "

$ws.Range("E38").Value = "  +2.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0429"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.128"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.94"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").Value = "0.0₃0688"
$ws.Range("E43").Value = "  +16.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.304"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.98%  "

$ws.Range("D46").Value = "2.915.69"
$ws.Range("E46").Value = "  -2.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.29%  "

$ws.Range("E48").Value = "  +11.20%  "

$ws.Range("E49").Value = "  +4.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.50%  "
